$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# need NumberFormat forced to Text ("@") first so they remain stored as strings,
# matching the original inline-string cell type in the workbook.

$ws.Range('D2').Value = '45.976.87'
$ws.Range('E2').Value = '  -1.50%  '

$ws.Range('D3').Value = '2.387.25'
$ws.Range('E3').Value = '  +3.40%  '

$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.66'
$ws.Range('E5').Value = '  -0.15%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.19'
$ws.Range('E6').Value = '  -2.92%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  -0.92%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -3.77%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.52'
$ws.Range('E10').Value = '  -6.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  -1.66%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.14'
$ws.Range('E12').Value = '  -3.36%  '

$ws.Range('E13').Value = '  -0.25%  '

$ws.Range('D14').Value = '2.747.75'
$ws.Range('E14').Value = '  +3.17%  '

$ws.Range('D15').Value = '2.384.80'
$ws.Range('E15').Value = '  +3.42%  '

$ws.Range('E16').Value = '  -0.35%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.72'
$ws.Range('E17').Value = '  -2.32%  '

$ws.Range('D18').Value = '45.900.50'
$ws.Range('E18').Value = '  -1.61%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.76'
$ws.Range('E19').Value = '  -3.51%  '

$ws.Range('D20').Value = '0.0₃0954'
$ws.Range('E20').Value = '  +0.75%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.05'
$ws.Range('E21').Value = '  -1.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.04'
$ws.Range('E22').Value = '  +0.27%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '244.05'
$ws.Range('E23').Value = '  -1.67%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.80'
$ws.Range('E24').Value = '  -5.14%  '

$ws.Range('E25').Value = '  +0.13%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.92'
$ws.Range('E26').Value = '  -2.06%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '39.65'
$ws.Range('E27').Value = '  -9.31%  '

$ws.Range('E28').Value = '  -2.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.76'
$ws.Range('E29').Value = '  -1.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.81'
$ws.Range('E30').Value = '  +20.74%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '21.11'
$ws.Range('E31').Value = '  +4.70%  '

$ws.Range('E32').Value = '  +7.10%  '

$ws.Range('E33').Value = '  -4.63%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '146.87'
$ws.Range('E34').Value = '  +0.73%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0772'
$ws.Range('E35').Value = '  -4.01%  '

$ws.Range('E36').Value = '  +0.74%  '

$ws.Range('E37').Value = '  +6.54%  '

$ws.Range('E38').Value = '  -3.33%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.84'
$ws.Range('E39').Value = '  -4.89%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.91'
$ws.Range('E40').Value = '  -3.89%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0299'
$ws.Range('E41').Value = '  -2.02%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.23'
$ws.Range('E42').Value = '  -6.92%  '

$ws.Range('D43').Value = '1.940.40'
$ws.Range('E43').Value = '  +4.68%  '

$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.49'
$ws.Range('E45').Value = '  +3.53%  '

$ws.Range('E46').Value = '  -9.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.49'
$ws.Range('E47').Value = '  +5.79%  '

$ws.Range('E48').Value = '  -5.13%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.66'
$ws.Range('E49').Value = '  +1.18%  '

$ws.Range('D50').Value = '2.618.45'
$ws.Range('E50').Value = '  +3.07%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '68.63'
$ws.Range('E51').Value = '  -8.29%  '
